$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-16 Friday", "2024-08-17 Saturday"),
    @("1+45=", "83-72="),
    @("72+11=", "62+13="),
    @("6+9=", "68-33="),
    @("97-31=", "31+12="),
    @("57-22=", "72-51="),
    @("76-60=", "49-4="),
    @("55-7=", "30+3="),
    @("90-54=", "0+29="),
    @("63+25=", "33+28="),
    @("10+77=", "9+67="),
    @("14+32=", "81-24="),
    @("56-50=", "77-71="),
    @("15+34=", "75-25="),
    @("52+0=", "35-4="),
    @("75-62=", "48-43="),
    @("10+86=", "96-5="),
    @("41+58=", "63-18="),
    @("95-10=", "70+17="),
    @("35+36=", "23+36="),
    @("46-35=", "42+55="),
    @("41-6=", "90-70="),
    @("90+2=", "31-8="),
    @("16-15=", "8+75="),
    @("62-46=", "3+75="),
    @("57+4=", "31+59="),
    @("98-46=", "66-56="),
    @("37-9=", "26+34="),
    @("16-6=", "59-44="),
    @("93-22=", "14-1="),
    @("53+29=", "5+25="),
    @("96-89=", "76+10="),
    @("74-5=", "42+4="),
    @("67+14=", "76+19="),
    @("99-17=", "28+36="),
    @("52-46=", "53+13="),
    @("93-30=", "75-62="),
    @("60+9=", "7+30="),
    @("28+14=", "19-3="),
    @("99-54=", "74-43="),
    @("96-67=", "62-32="),
    @("15-14=", "3+96="),
    @("36+4=", "43+48="),
    @("87-74=", "79-1="),
    @("35+5=", "37+59="),
    @("65-11=", "41-21="),
    @("27-21=", "97-69="),
    @("24+12=", "39+51="),
    @("82+12=", "35+49="),
    @("72+1=", "89-75="),
    @("13+39=", "84-75="),
    @("45+45=", "73-66="),
    @("16+19=", "2+73="),
    @("37-30=", "90-70="),
    @("53-15=", "93-80="),
    @("78-61=", "14+36="),
    @("49+1=", "6+35="),
    @("16+58=", "9+68="),
    @("70-49=", "78+20="),
    @("11+72=", "89-46="),
    @("36+23=", "19+32="),
    @("67-55=", "53+39="),
    @("20+13=", "61+36="),
    @("85+7=", "49-13="),
    @("69-48=", "80-48="),
    @("92-80=", "61-12="),
    @("82-4=", "14+64="),
    @("22+55=", "15+63="),
    @("94-35=", "22+5="),
    @("73-53=", "81-36="),
    @("52+22=", "21-13="),
    @("13+20=", "10+32="),
    @("32+8=", "89+4="),
    @("74-52=", "17+35="),
    @("17-3=", "99-36="),
    @("2+11=", "64-50="),
    @("77+13=", "71+28="),
    @("39+60=", "97-83="),
    @("69+24=", "58-44="),
    @("88-1=", "83-31="),
    @("75+9=", "9+61="),
    @("26+4=", "75-66="),
    @("20+18=", "61-53="),
    @("74+10=", "67-11="),
    @("88-48=", "30+7="),
    @("75-50=", "88-6="),
    @("48-34=", "20+16="),
    @("77-30=", "70+3="),
    @("85+2=", "56-33="),
    @("0+53=", "71-67="),
    @("85-40=", "5+12="),
    @("97-0=", "17+23="),
    @("58+20=", "12+50="),
    @("7+5=", "25+53="),
    @("58-18=", "22-19="),
    @("94-3=", "45+23="),
    @("47+40=", "61-34="),
    @("25+37=", "57-17="),
    @("73-64=", "33+20="),
    @("31+4=", "46-41="),
    @("61-0=", "13+18="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
